# "new format for excel"
#
# The sheet's Public/Private/Save/Cache/Ref/Desc rows (rows 3-9, columns
# B:G) were stored as Boolean cells (t="b", 0/1 meaning FALSE/TRUE). The
# new format stores the same 0 values as plain numbers instead, so Excel
# (and downstream tooling) reads them back as numeric rather than boolean.
#
# Re-assigning via Value2 re-types the cell: writing a plain number clears
# the boolean flag while leaving the literal value (0) and the existing
# cell style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:G9").Value2 = 0
